# Apply the edits described by the diff to the active document.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...que acontece em Salvador (BA), até sexta-feira (23)."
#    -> "...que aconteceu em Salvador (BA)."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "que acontece em Salvador (BA), até sexta-feira (23).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "que aconteceu em Salvador (BA).", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) "Carente em trabalhos que de fato" -> "Carente em realizações que de fato"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Carente em trabalhos que de fato",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Carente em realizações que de fato", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) "Por mais que melhorias venham" -> "Por mais, que as melhorias venham"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Por mais que melhorias venham",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Por mais, que as melhorias venham", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) "vem sendo proposto." -> "vem sendo proposto atualmente."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "vem sendo proposto.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "vem sendo proposto atualmente.", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Insert a brand-new paragraph right after the "Carente em ..." one
#    (i.e. right before the blank paragraph that precedes "Pensando ...").
# ---------------------------------------------------------------------
$find = $d.Content
$find.Find.Execute("Carente em realizações", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$carenteParaIndex = $find.Paragraphs(1).Index

$carentePara = $d.Paragraphs($carenteParaIndex)
$carentePara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($carenteParaIndex + 1)

$newParaXml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:jc w:val='both'/>
    <w:rPr>
      <w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>
      <w:bCs/>
      <w:color w:val='000000' w:themeColor='text1'/>
      <w:sz w:val='24'/>
      <w:szCs w:val='24'/>
      <w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/>
      <w:bCs/>
      <w:color w:val='000000' w:themeColor='text1'/>
      <w:sz w:val='24'/>
      <w:szCs w:val='24'/>
      <w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t xml:space='preserve'>Importante deixar claro o quanto pesquisas relacionadas a estes temas podem ajudar tanto na prática quanto na teoria. Empiricamente surgiram problemas, inevitavelmente e esses problemas terão de ser solucionados para assim alcançar um melhor resultado para a diminuição do desperdício alimentar. Desse modo, entra a parte teórica que deverá solucionar de alguma maneira esses problemas com novas tecnologias, tomadas de decisão, etc...</w:t>
  </w:r>
</w:p>
"@

$newPara.Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------
# 6) Move the "_GoBack" bookmark from the end of the "Carente em ..."
#    paragraph to the end of the freshly inserted paragraph.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$newPara2 = $d.Paragraphs($carenteParaIndex + 1)
$bmStart = $newPara2.Range.End - 1
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------
# 7) Remove the lastRenderedPageBreak that used to sit in front of
#    "Pensando " (it now belongs on the new paragraph above instead).
# ---------------------------------------------------------------------
$find2 = $d.Content
$ok = $find2.Find.Execute("Pensando ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok) {
    $find2.Delete()
    $find2.InsertAfter("Pensando ")
}

Write-Output "done"
